$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (C) column date from 45175 (2023-09-06) to 45183 (2023-09-14)
# for rows 2 through 6, preserving existing number formatting/style.
foreach ($r in 2..6) {
    $ws.Range("C$r").Value = 45183
}
